$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append (CME ETHUSD RR raw data refresh)
$data = @(
    @("2025-09-06", 4297.78, 4266.18, 4306.97),
    @("2025-09-07", 4303.15, 4286.82, 4292.52),
    @("2025-09-08", 4337.28, 4303.83, 4297.15),
    @("2025-09-09", 4311.4,  4292.99, 4367.05),
    @("2025-09-10", 4427.17, 4320.38, 4323.38),
    @("2025-09-11", 4422.02, 4422.98, 4440.36),
    @("2025-09-12", 4549.6,  4638.02, 4538.11)
)

$startRow = 394
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    $dateCell = $ws.Cells.Item($row, 1)
    # Force text format so the date-like string isn't auto-converted to a
    # serial date, then drop back to the default "Normal" style so the
    # cell doesn't end up carrying an explicit style index, matching the
    # plain text cells already used throughout column A.
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $rec[0]
    $dateCell.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
}
